{"js": "// Replace the date line and the 25 division-problem answers in the table.\n// Each \"before\" text is unique in the document, so a simple exact search\n// and replace (matchCase + matchWholeWord) is sufficient and unambiguous.\nconst replacements = [\n  [\"2025-03-01 Saturday\", \"2025-03-02 Sunday\"],\n  [\"96\u00f74=24, 0\", \"91\u00f72=45, 1\"],\n  [\"49\u00f79=5, 4\", \"52\u00f75=10, 2\"],\n  [\"38\u00f75=7, 3\", \"89\u00f77=12, 5\"],\n  [\"67\u00f72=33, 1\", \"60\u00f72=30, 0\"],\n  [\"82\u00f72=41, 0\", \"34\u00f75=6, 4\"],\n  [\"15\u00f79=1, 6\", \"77\u00f77=11, 0\"],\n  [\"78\u00f77=11, 1\", \"68\u00f74=17, 0\"],\n  [\"13\u00f72=6, 1\", \"71\u00f77=10, 1\"],\n  [\"63\u00f78=7, 7\", \"30\u00f79=3, 3\"],\n  [\"46\u00f76=7, 4\", \"96\u00f72=48, 0\"],\n  [\"83\u00f77=11, 6\", \"50\u00f74=12, 2\"],\n  [\"17\u00f74=4, 1\", \"97\u00f76=16, 1\"],\n  [\"65\u00f74=16, 1\", \"51\u00f76=8, 3\"],\n  [\"44\u00f79=4, 8\", \"58\u00f73=19, 1\"],\n  [\"51\u00f72=25, 1\", \"32\u00f78=4, 0\"],\n  [\"65\u00f73=21, 2\", \"25\u00f75=5, 0\"],\n  [\"32\u00f73=10, 2\", \"24\u00f72=12, 0\"],\n  [\"65\u00f72=32, 1\", \"31\u00f77=4, 3\"],\n  [\"36\u00f73=12, 0\", \"71\u00f79=7, 8\"],\n  [\"59\u00f78=7, 3\", \"11\u00f78=1, 3\"],\n  [\"46\u00f75=9, 1\", \"22\u00f78=2, 6\"],\n  [\"80\u00f78=10, 0\", \"28\u00f72=14, 0\"],\n  [\"21\u00f76=3, 3\", \"94\u00f72=47, 0\"],\n  [\"86\u00f79=9, 5\", \"21\u00f78=2, 5\"],\n  [\"22\u00f73=7, 1\", \"62\u00f78=7, 6\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${before}\"`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 division-problem answers in the table.\n# Each \"before\" text is unique in the document, so Find/Replace with\n# MatchWholeWord (exact text) and ReplaceAll is unambiguous.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-03-01 Saturday\", \"2025-03-02 Sunday\"),\n    @(\"96\u00f74=24, 0\", \"91\u00f72=45, 1\"),\n    @(\"49\u00f79=5, 4\", \"52\u00f75=10, 2\"),\n    @(\"38\u00f75=7, 3\", \"89\u00f77=12, 5\"),\n    @(\"67\u00f72=33, 1\", \"60\u00f72=30, 0\"),\n    @(\"82\u00f72=41, 0\", \"34\u00f75=6, 4\"),\n    @(\"15\u00f79=1, 6\", \"77\u00f77=11, 0\"),\n    @(\"78\u00f77=11, 1\", \"68\u00f74=17, 0\"),\n    @(\"13\u00f72=6, 1\", \"71\u00f77=10, 1\"),\n    @(\"63\u00f78=7, 7\", \"30\u00f79=3, 3\"),\n    @(\"46\u00f76=7, 4\", \"96\u00f72=48, 0\"),\n    @(\"83\u00f77=11, 6\", \"50\u00f74=12, 2\"),\n    @(\"17\u00f74=4, 1\", \"97\u00f76=16, 1\"),\n    @(\"65\u00f74=16, 1\", \"51\u00f76=8, 3\"),\n    @(\"44\u00f79=4, 8\", \"58\u00f73=19, 1\"),\n    @(\"51\u00f72=25, 1\", \"32\u00f78=4, 0\"),\n    @(\"65\u00f73=21, 2\", \"25\u00f75=5, 0\"),\n    @(\"32\u00f73=10, 2\", \"24\u00f72=12, 0\"),\n    @(\"65\u00f72=32, 1\", \"31\u00f77=4, 3\"),\n    @(\"36\u00f73=12, 0\", \"71\u00f79=7, 8\"),\n    @(\"59\u00f78=7, 3\", \"11\u00f78=1, 3\"),\n    @(\"46\u00f75=9, 1\", \"22\u00f78=2, 6\"),\n    @(\"80\u00f78=10, 0\", \"28\u00f72=14, 0\"),\n    @(\"21\u00f76=3, 3\", \"94\u00f72=47, 0\"),\n    @(\"86\u00f79=9, 5\", \"21\u00f78=2, 5\"),\n    @(\"22\u00f73=7, 1\", \"62\u00f78=7, 6\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $pair[0]\n    $replace = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.MatchWildcards = $false\n    $range.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
